$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'reizen',
    'scheitern',
    'wechseln',
    'spüren',
    'wandern',
    'schwören',
    'enden',
    'merken',
    'treiben',
    'rasen',
    'wachsen',
    'nennen',
    'seufzen',
    'töten',
    'graben',
    'drücken',
    'dringen',
    'fällen',
    'herzen',
    'tropfen',
    'räumen',
    'rufen',
    'gründen',
    'ehren',
    'schmecken',
    'wiegen',
    'decken',
    'kehren',
    'stoppen',
    'tollen',
    'lassen',
    'sorgen',
    'schwimmen',
    'baden',
    'zielen',
    'sterben',
    'lösen',
    'kümmern',
    'saufen',
    'nerven',
    'trennen',
    'herrschen',
    'stören',
    'betteln',
    'klagen',
    'machen',
    'boxen',
    'zeugen',
    'flehen',
    'mögen',
    'lehnen',
    'zahlen',
    'wenden',
    'fragen',
    'erben',
    'heben',
    'stecken',
    'schneiden',
    'schämen',
    'laufen',
    'bluten',
    'trauen',
    'streichen',
    'wohnen',
    'dienen',
    'gelten',
    'bieten',
    'bellen',
    'sperren',
    'danken',
    'schenken',
    'wehtun',
    'süßen',
    'biegen',
    'fahren',
    'backen',
    'fesseln',
    'knarren',
    'klingen',
    'flüchten',
    'kürzen',
    'fließen',
    'weichen',
    'retten',
    'heilen',
    'pfeifen',
    'münzen',
    'schätzen',
    'schrecken',
    'irren',
    'formen',
    'buchen',
    'kosten',
    'schwächen',
    'mauern',
    'regeln',
    'jubeln',
    'sprengen',
    'stürmen',
    'sichern',
    'stoßen',
    'hören',
    'zögern',
    'feiern',
    'geben',
    'wehren',
    'testen',
    'stillen',
    'fangen',
    'ärgern',
    'helfen',
    'bitten',
    'zeigen',
    'arten',
    'schleppen',
    'filmen',
    'bergen',
    'suchen',
    'warnen',
    'wellen',
    'kichern',
    'altern',
    'stehlen',
    'ändern',
    'lieben',
    'siegen',
    'starten',
    'brauchen',
    'werden',
    'rühren',
    'sammeln',
    'bilden',
    'malen',
    'schlucken',
    'holen',
    'drehen',
    'planen',
    'lügen',
    'reisen',
    'runden',
    'führen',
    'jagen',
    'schütteln',
    'wüten',
    'fallen',
    'pflanzen',
    'quälen',
    'platzen',
    'rauchen',
    'liegen',
    'teilen',
    'sinken',
    'liefern',
    'wirken',
    'knurren',
    'klettern',
    'trösten',
    'heulen',
    'grüßen',
    'werfen',
    'parken',
    'scheinen',
    'fordern',
    'lockern',
    'schulden',
    'achten',
    'flüstern',
    'haben',
    'folgen',
    'spinnen',
    'stammen',
    'foltern',
    'loben',
    'lesen',
    'schwingen',
    'äußern',
    'beißen',
    'spielen',
    'wundern',
    'bauen',
    'fischen',
    'freuen',
    'pressen',
    'zünden',
    'schreiten',
    'stürzen',
    'nutzen',
    'greifen',
    'gleichen',
    'rauben',
    'hauen',
    'opfern'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}
